# BOT; UPDATE DATA
# Appends one new daily row (2020-05-23, serial 43974) to the three data
# sheets ("all", "kobe", "other") just above their trailing footnote row.
# Inserting the row (rather than just overwriting the footnote row and
# retyping the footnote one row down) makes Excel copy the formatting of
# the row above for the new row and push the footnote row down by one -
# exactly what the diff shows.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet "all": new data lives at row 46 (old footnote row 46 -> 47)
# ---------------------------------------------------------------------------
$wsAll = $wb.Worksheets.Item("all")
$wsAll.Activate() | Out-Null
$wsAll.Rows.Item(46).Insert()
$wsAll.Cells.Item(46, 1).Value = 43974
$wsAll.Cells.Item(46, 2).Value = 285
$wsAll.Cells.Item(46, 3).Value = 282
$wsAll.Cells.Item(46, 4).Value = 36
$wsAll.Cells.Item(46, 5).Value = 32
$wsAll.Cells.Item(46, 6).Value = 4
$wsAll.Cells.Item(46, 7).Value = 12
$wsAll.Cells.Item(46, 8).Value = 234
$wsAll.Range("C49").Select() | Out-Null

# ---------------------------------------------------------------------------
# Sheet "kobe": new data lives at row 101 (old footnote row 101 -> 102)
# ---------------------------------------------------------------------------
$wsKobe = $wb.Worksheets.Item("kobe")
$wsKobe.Activate() | Out-Null
$wsKobe.Rows.Item(101).Insert()
$wsKobe.Cells.Item(101, 1).Value = 43974
$wsKobe.Cells.Item(101, 2).Value = 0
$wsKobe.Cells.Item(101, 3).Value = 3010
$wsKobe.Cells.Item(101, 4).Value = 0
$wsKobe.Cells.Item(101, 5).Value = 285
$wsKobe.Cells.Item(101, 6).Value = 31
$wsKobe.Cells.Item(101, 7).Value = 28
$wsKobe.Cells.Item(101, 8).Value = 3
$wsKobe.Cells.Item(101, 9).Value = 12
$wsKobe.Cells.Item(101, 10).Value = 225
$wsKobe.Range("K101").Select() | Out-Null

# ---------------------------------------------------------------------------
# Sheet "other": new data lives at row 76 (old footnote row 76 -> 77)
# ---------------------------------------------------------------------------
$wsOther = $wb.Worksheets.Item("other")
$wsOther.Activate() | Out-Null
$wsOther.Rows.Item(76).Insert()
$wsOther.Cells.Item(76, 1).Value = 43974
$wsOther.Cells.Item(76, 2).Value = 0
$wsOther.Cells.Item(76, 3).Value = 14
$wsOther.Cells.Item(76, 4).Value = 5
$wsOther.Cells.Item(76, 5).Value = 4
$wsOther.Cells.Item(76, 6).Value = 1
$wsOther.Cells.Item(76, 7).Value = 0
$wsOther.Cells.Item(76, 8).Value = 9
$wsOther.Range("F77").Select() | Out-Null

# ---------------------------------------------------------------------------
# Restore "all" as the active/selected sheet (tabSelected="1"), as in the
# diff - it must remain the tab that was active before the edit.
# ---------------------------------------------------------------------------
$wsAll.Activate() | Out-Null
